$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerate the "K" column (column G) values - was previously based on Strike#,
# now recalculated. Update each row's G cell with the new K value.
$ws.Range("G2").Value = 0
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 0
$ws.Range("G8").Value = 2
$ws.Range("G9").Value = 2
$ws.Range("G10").Value = 3
$ws.Range("G12").Value = 1
